# "added 4wk low sales check"
# Updates forecast values (MyForecast / Seasonality Index) on the
# "Forecast Comparison" sheet and the corresponding rollup totals on the
# "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison" ---------------------------------------
# Column D = MyForecast, Column L = Seasonality Index

$ws1.Range("L2").Value = 0.9

$ws1.Range("D3").Value = 13
$ws1.Range("L3").Value = 1.05

$ws1.Range("D4").Value = 14
$ws1.Range("L4").Value = 1.18

$ws1.Range("D5").Value = 14
$ws1.Range("L5").Value = 0.8100000000000001

$ws1.Range("D6").Value = 13
$ws1.Range("L6").Value = 1.04

$ws1.Range("D7").Value = 11
$ws1.Range("L7").Value = 1.14

$ws1.Range("L8").Value = 1

$ws1.Range("D9").Value = 13
$ws1.Range("L9").Value = 1.1

$ws1.Range("D10").Value = 14
$ws1.Range("L10").Value = 1.11

$ws1.Range("D11").Value = 14
$ws1.Range("L11").Value = 0.95

$ws1.Range("L12").Value = 1.15

$ws1.Range("D13").Value = 10
$ws1.Range("L13").Value = 0.83

$ws1.Range("D14").Value = 11
$ws1.Range("L14").Value = 0.91

$ws1.Range("D15").Value = 13
$ws1.Range("L15").Value = 0.95

$ws1.Range("D16").Value = 14
$ws1.Range("L16").Value = 1.2

$ws1.Range("D17").Value = 14
$ws1.Range("L17").Value = 1.12

# --- Sheet "Summary" -----------------------------------------------------
# Column B values are stored as text; force text formatting before
# assigning numeric-looking strings so they remain text cells.

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "211"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "106"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "55"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "14"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "11"
